$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# Weekly crime statistics table updates (rows 14-29)
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -66.666666666666
$ws.Range("J14").Value = 14
$ws.Range("K14").Value = -28.571428571428

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 35.294117647058
$ws.Range("L15").Value = 130
$ws.Range("M15").Value = 35.294117647058
$ws.Range("N15").Value = -47.727272727272

# Row 16
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 47
$ws.Range("G16").Value = 51
$ws.Range("H16").Value = -7.843137254901
$ws.Range("I16").Value = 285
$ws.Range("J16").Value = 323
$ws.Range("K16").Value = -11.764705882352
$ws.Range("L16").Value = 41.791044776119
$ws.Range("M16").Value = 22.844827586206
$ws.Range("N16").Value = -71.669980119284

# Row 17
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = -4
$ws.Range("F17").Value = 93
$ws.Range("G17").Value = 108
$ws.Range("H17").Value = -13.888888888888
$ws.Range("I17").Value = 546
$ws.Range("J17").Value = 511
$ws.Range("K17").Value = 6.849315068493
$ws.Range("L17").Value = 39.641943734015
$ws.Range("M17").Value = 98.545454545454
$ws.Range("N17").Value = -18.018018018018

# Row 18
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 700
$ws.Range("F18").Value = 36
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 260
$ws.Range("J18").Value = 209
$ws.Range("K18").Value = 24.401913875598
$ws.Range("L18").Value = 106.349206349206
$ws.Range("M18").Value = 104.724409448819
$ws.Range("N18").Value = -74.708171206225

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = -36.904761904761
$ws.Range("I19").Value = 366
$ws.Range("J19").Value = 457
$ws.Range("K19").Value = -19.912472647702
$ws.Range("L19").Value = 5.475504322766
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = -8.5

# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 38
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = 11.764705882352
$ws.Range("I20").Value = 246
$ws.Range("J20").Value = 217
$ws.Range("K20").Value = 13.364055299539
$ws.Range("L20").Value = 267.164179104478
$ws.Range("M20").Value = 143.564356435644
$ws.Range("N20").Value = -65.975103734439

# Row 21
$ws.Range("C21").Value = 71
$ws.Range("D21").Value = 73
$ws.Range("E21").Value = -2.739726027397
$ws.Range("F21").Value = 272
$ws.Range("G21").Value = 306
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 1736
$ws.Range("J21").Value = 1748
$ws.Range("K21").Value = -0.686498855835
$ws.Range("L21").Value = 50.694444444444
$ws.Range("M21").Value = 84.093319194061
$ws.Range("N21").Value = -55.498590105101

# Row 22
$ws.Range("D22").Value = 4
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -62.5
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 9.523809523809
$ws.Range("M22").Value = 0

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -54.545454545454
$ws.Range("I23").Value = 38
$ws.Range("J23").Value = 37
$ws.Range("K23").Value = 2.702702702702
$ws.Range("L23").Value = 5.555555555555
$ws.Range("M23").Value = 35.714285714285

# Row 24
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 145
$ws.Range("H24").Value = 28.965517241379
$ws.Range("I24").Value = 1063
$ws.Range("J24").Value = 934
$ws.Range("K24").Value = 13.811563169164
$ws.Range("L24").Value = 57.481481481481
$ws.Range("M24").Value = 55.865102639296

# Row 25
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = -29.032258064516
$ws.Range("F25").Value = 103
$ws.Range("G25").Value = 120
$ws.Range("H25").Value = -14.166666666666
$ws.Range("I25").Value = 621
$ws.Range("J25").Value = 651
$ws.Range("K25").Value = -4.608294930875
$ws.Range("L25").Value = 28.838174273858
$ws.Range("M25").Value = 1.305057096247

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 29
$ws.Range("K26").Value = 24.137931034482
$ws.Range("L26").Value = 63.636363636363

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 72
$ws.Range("J27").Value = 56
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 60

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 31
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = -11.428571428571
$ws.Range("L28").Value = -27.906976744186
$ws.Range("M28").Value = 63.157894736842
$ws.Range("N28").Value = -64.772727272727

# Row 29
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 26
$ws.Range("J29").Value = 31
$ws.Range("K29").Value = -16.129032258064
$ws.Range("L29").Value = -23.529411764705
$ws.Range("M29").Value = 52.941176470588
$ws.Range("N29").Value = -66.666666666666
